# OPTIMIZACION DE INTERFAZ NO.4
# Adds 13 new login/logout records (rows 107-119) to the "Logins" sheet,
# mirroring the formatting of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of (Tiempo Entrada, Tiempo Salida) for the new rows 107-119.
$entries = @(
    @("20/10/2024 17:43:05", "20/10/2024 17:43:12"),
    @("20/10/2024 10:31:19", "20/10/2024 10:32:42"),
    @("20/10/2024 14:14:55", "20/10/2024 14:14:58"),
    @("20/10/2024 14:37:00", "20/10/2024 14:40:07"),
    @("20/10/2024 14:40:24", "20/10/2024 14:41:48"),
    @("20/10/2024 14:42:04", "20/10/2024 14:49:01"),
    @("20/10/2024 14:53:44", "20/10/2024 14:54:10"),
    @("20/10/2024 14:57:54", "20/10/2024 15:49:58"),
    @("20/10/2024 14:58:00", "20/10/2024 17:33:23"),
    @("20/10/2024 14:58:01", "20/10/2024 17:35:56"),
    @("20/10/2024 17:42:29", "20/10/2024 17:42:34"),
    @("20/10/2024 17:42:43", "20/10/2024 17:42:55"),
    @("20/10/2024 17:43:05", "20/10/2024 17:43:12")
)

$personal = "javiergonzalezcoradopineed"
$rol = "ADMINISTRADOR"

$startRow = 107
for ($i = 0; $i -lt $entries.Count; $i++) {
    $row = $startRow + $i
    $entrada = $entries[$i][0]
    $salida = $entries[$i][1]

    # Copy formatting (style/number format) from an existing data row (row 2)
    # into columns A and B of the new row, then overwrite the value with the
    # new text so it keeps the same style index as other entries.
    $ws.Range("A2").Copy($ws.Range("A$row"))
    $ws.Range("A$row").Value = $entrada

    $ws.Range("B2").Copy($ws.Range("B$row"))
    $ws.Range("B$row").Value = $salida

    $ws.Range("C$row").Value = $personal
    $ws.Range("D$row").Value = $rol
}

# Update the view to reflect scrolling down to the newly added rows.
$ws.Activate()
$ws.Range("C106:C119").Select()
$excel.ActiveWindow.ScrollRow = 107
